$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("S2").Value = 1.75
$ws.Range("T2").Value = 2.05

# Row 3 updates
$ws.Range("G3").Value = 1.53
$ws.Range("I3").Value = 5.25
$ws.Range("J3").Value = 2.05
$ws.Range("W3").Value = 9.5
$ws.Range("Z3").Value = 12
$ws.Range("AI3").Value = 29
$ws.Range("AM3").Value = 34
$ws.Range("AV3").Value = 41
$ws.Range("AX3").Value = 26
$ws.Range("AZ3").Value = 81

$wb.Save()
